$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Mmp2'
$ws.Cells.Item(2, 3).Value = 'Sdc2'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.332428
$ws.Cells.Item(2, 8).Value = 9.997284000000001
$ws.Cells.Item(2, 9).Value = 0.01078284025505985
$ws.Cells.Item(2, 10).Value = 0.01078284025505985
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 1.321445333333333
$ws.Cells.Item(2, 14).Value = 3.964336
$ws.Cells.Item(2, 15).Value = 0.01021782062667047
$ws.Cells.Item(2, 16).Value = 0.01021782062667047
$ws.Cells.Item(2, 17).Value = 4.403621429269334
$ws.Cells.Item(2, 18).Value = 39.63259286342401
$ws.Cells.Item(2, 19).Value = 0.0001101771275722432
$ws.Cells.Item(2, 20).Value = 0.0001101771275722432

# Row 3
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Mmp2'
$ws.Cells.Item(3, 3).Value = 'Sdc2'
$ws.Cells.Item(3, 4).Value = 'FAPs'
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.332428
$ws.Cells.Item(3, 8).Value = 9.997284000000001
$ws.Cells.Item(3, 9).Value = 0.01078284025505985
$ws.Cells.Item(3, 10).Value = 0.01078284025505985
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 105.9632263333333
$ws.Cells.Item(3, 14).Value = 317.889679
$ws.Cells.Item(3, 15).Value = 0.819340166699254
$ws.Cells.Item(3, 16).Value = 0.8193401666992541
$ws.Cells.Item(3, 17).Value = 353.1148224035373
$ws.Cells.Item(3, 18).Value = 3178.033401631836
$ws.Cells.Item(3, 19).Value = 0.008834814132072167
$ws.Cells.Item(3, 20).Value = 0.008834814132072169

# Row 4
$ws.Cells.Item(4, 1).Value = 'ECs'
$ws.Cells.Item(4, 2).Value = 'Mmp2'
$ws.Cells.Item(4, 3).Value = 'Sdc2'
$ws.Cells.Item(4, 4).Value = 'sCs'
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.332428
$ws.Cells.Item(4, 8).Value = 9.997284000000001
$ws.Cells.Item(4, 9).Value = 0.01078284025505985
$ws.Cells.Item(4, 10).Value = 0.01078284025505985
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 22.04284166666666
$ws.Cells.Item(4, 14).Value = 66.128525
$ws.Cells.Item(4, 15).Value = 0.1704420126740755
$ws.Cells.Item(4, 16).Value = 0.1704420126740755
$ws.Cells.Item(4, 17).Value = 73.45618276956667
$ws.Cells.Item(4, 18).Value = 661.1056449261
$ws.Cells.Item(4, 19).Value = 0.001837848995415443
$ws.Cells.Item(4, 20).Value = 0.001837848995415443

# Row 5
$ws.Cells.Item(5, 1).Value = 'FAPs'
$ws.Cells.Item(5, 2).Value = 'Mmp2'
$ws.Cells.Item(5, 3).Value = 'Sdc2'
$ws.Cells.Item(5, 4).Value = 'ECs'
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 295.7980143333334
$ws.Cells.Item(5, 8).Value = 887.394043
$ws.Cells.Item(5, 9).Value = 0.9571227754418815
$ws.Cells.Item(5, 10).Value = 0.9571227754418815
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 1.321445333333333
$ws.Cells.Item(5, 14).Value = 3.964336
$ws.Cells.Item(5, 15).Value = 0.01021782062667047
$ws.Cells.Item(5, 16).Value = 0.01021782062667047
$ws.Cells.Item(5, 17).Value = 390.8809056500498
$ws.Cells.Item(5, 18).Value = 3517.928150850448
$ws.Cells.Item(5, 19).Value = 0.009779708837166144
$ws.Cells.Item(5, 20).Value = 0.009779708837166145

# Row 6
$ws.Cells.Item(6, 1).Value = 'FAPs'
$ws.Cells.Item(6, 2).Value = 'Mmp2'
$ws.Cells.Item(6, 3).Value = 'Sdc2'
$ws.Cells.Item(6, 4).Value = 'FAPs'
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 295.7980143333334
$ws.Cells.Item(6, 8).Value = 887.394043
$ws.Cells.Item(6, 9).Value = 0.9571227754418815
$ws.Cells.Item(6, 10).Value = 0.9571227754418815
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 105.9632263333333
$ws.Cells.Item(6, 14).Value = 317.889679
$ws.Cells.Item(6, 15).Value = 0.819340166699254
$ws.Cells.Item(6, 16).Value = 0.8193401666992541
$ws.Cells.Item(6, 17).Value = 31343.71194175358
$ws.Cells.Item(6, 18).Value = 282093.4074757822
$ws.Cells.Item(6, 19).Value = 0.7842091343822039
$ws.Cells.Item(6, 20).Value = 0.784209134382204

# Row 7
$ws.Cells.Item(7, 1).Value = 'FAPs'
$ws.Cells.Item(7, 2).Value = 'Mmp2'
$ws.Cells.Item(7, 3).Value = 'Sdc2'
$ws.Cells.Item(7, 4).Value = 'sCs'
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 295.7980143333334
$ws.Cells.Item(7, 8).Value = 887.394043
$ws.Cells.Item(7, 9).Value = 0.9571227754418815
$ws.Cells.Item(7, 10).Value = 0.9571227754418815
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 22.04284166666666
$ws.Cells.Item(7, 14).Value = 66.128525
$ws.Cells.Item(7, 15).Value = 0.1704420126740755
$ws.Cells.Item(7, 16).Value = 0.1704420126740755
$ws.Cells.Item(7, 17).Value = 6520.228795264064
$ws.Cells.Item(7, 18).Value = 58682.05915737657
$ws.Cells.Item(7, 19).Value = 0.1631339322225115
$ws.Cells.Item(7, 20).Value = 0.1631339322225115

# Row 8
$ws.Cells.Item(8, 1).Value = 'sCs'
$ws.Cells.Item(8, 2).Value = 'Mmp2'
$ws.Cells.Item(8, 3).Value = 'Sdc2'
$ws.Cells.Item(8, 4).Value = 'ECs'
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 9.918743333333333
$ws.Cells.Item(8, 8).Value = 29.75623
$ws.Cells.Item(8, 9).Value = 0.03209438430305867
$ws.Cells.Item(8, 10).Value = 0.03209438430305867
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 1.321445333333333
$ws.Cells.Item(8, 14).Value = 3.964336
$ws.Cells.Item(8, 15).Value = 0.01021782062667047
$ws.Cells.Item(8, 16).Value = 0.01021782062667047
$ws.Cells.Item(8, 17).Value = 13.10707709036445
$ws.Cells.Item(8, 18).Value = 117.96369381328
$ws.Cells.Item(8, 19).Value = 0.0003279346619320818
$ws.Cells.Item(8, 20).Value = 0.0003279346619320818

# Row 9
$ws.Cells.Item(9, 1).Value = 'sCs'
$ws.Cells.Item(9, 2).Value = 'Mmp2'
$ws.Cells.Item(9, 3).Value = 'Sdc2'
$ws.Cells.Item(9, 4).Value = 'FAPs'
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 9.918743333333333
$ws.Cells.Item(9, 8).Value = 29.75623
$ws.Cells.Item(9, 9).Value = 0.03209438430305867
$ws.Cells.Item(9, 10).Value = 0.03209438430305867
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 105.9632263333333
$ws.Cells.Item(9, 14).Value = 317.889679
$ws.Cells.Item(9, 15).Value = 0.819340166699254
$ws.Cells.Item(9, 16).Value = 0.8193401666992541
$ws.Cells.Item(9, 17).Value = 1051.022044772241
$ws.Cells.Item(9, 18).Value = 9459.198402950169
$ws.Cells.Item(9, 19).Value = 0.02629621818497802
$ws.Cells.Item(9, 20).Value = 0.02629621818497802

# Row 10
$ws.Cells.Item(10, 1).Value = 'sCs'
$ws.Cells.Item(10, 2).Value = 'Mmp2'
$ws.Cells.Item(10, 3).Value = 'Sdc2'
$ws.Cells.Item(10, 4).Value = 'sCs'
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 9.918743333333333
$ws.Cells.Item(10, 8).Value = 29.75623
$ws.Cells.Item(10, 9).Value = 0.03209438430305867
$ws.Cells.Item(10, 10).Value = 0.03209438430305867
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 22.04284166666666
$ws.Cells.Item(10, 14).Value = 66.128525
$ws.Cells.Item(10, 15).Value = 0.1704420126740755
$ws.Cells.Item(10, 16).Value = 0.1704420126740755
$ws.Cells.Item(10, 17).Value = 218.6372888289722
$ws.Cells.Item(10, 18).Value = 1967.73559946075
$ws.Cells.Item(10, 19).Value = 0.005470231456148575
$ws.Cells.Item(10, 20).Value = 0.005470231456148575
